$wb = $excel.ActiveWorkbook

# Sheet "展览" and "全部类型" both contain the same event listing data.
# Update the "想去人数" (number of people interested) counts:
#   Row 2 (丽水·CCAC动漫游戏嘉年华): 258 -> 259
#   Row 5 (丽水·动漫游戏展): 9 -> 10

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 259
    $ws.Range("F5").Value = 10
}
